$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values per row (B, C, D, E, F, G)
$data = @{
    2  = @(0.3342080079727115, 0.8001272719199065, 1.057948689255649, 1.02856632710567, 0.9994115000724915, 19)
    3  = @(0.2715251256614695, 0.5749987921201453, 0.6574421919879687, 0.8108280902805284, 0.7861631834124684, 18)
    4  = @(0.1126313932437445, 0.4526599150000412, 0.4170638627520434, 0.645804817845178, 0.6554781911006021, 17)
    5  = @(0.2541178104058748, 0.4548235677275672, 0.357578812029527, 0.5979789394531608, 0.5590496225693991, 16)
    6  = @(0.3625344086844117, 0.3625344086844117, 0.1936168257130608, 0.4400191197130652, 0.2581229529138574, 15)
    7  = @(0.3044989716643398, 0.3044989716643398, 0.1348154283086344, 0.3671722052506622, 0.212917691553823, 14)
    8  = @(0.3319890590706304, 0.3319890590706304, 0.1468494055944594, 0.3832093495655597, 0.1992119292935162, 13)
    9  = @(0.333398904683986, 0.3399377628279546, 0.1517941216966595, 0.3896076509729494, 0.2105558670468018, 12)
    10 = @(0.3614593346982711, 0.3614593346982711, 0.1605564165475934, 0.4006949170473633, 0.1813668175212442, 11)
    11 = @(0.3457552520087644, 0.3457552520087644, 0.1504834501066176, 0.3879219639394211, 0.1854027322553479, 10)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
